# "check in to run 4 scenarios"
# Clears the stale "Pass"/"Fail" status markers left over in column H
# (and the stray F7:F9 markers on the BatchDecisionOutputValidations
# checklist rows) so the next run starts from a clean slate, then moves
# the active/selected tab from BatchDecisionOutputValidations over to
# Cases_RealTimeSpine.

$wb = $excel.ActiveWorkbook

# --- VerifyCSVForNewVersion: clear the old pass/fail results in H5:H61 ---
$wsNewVersion = $wb.Worksheets.Item("VerifyCSVForNewVersion")
$wsNewVersion.Range("H5:H61").ClearContents()

# --- Cases_RealTimeSpine: clear the old pass/fail results in H43:H49 ---
$wsCases = $wb.Worksheets.Item("Cases_RealTimeSpine")
$wsCases.Range("H43:H49").ClearContents()

# --- BatchDecisionOutputValidations: drop the leftover F7:F9 markers ---
$wsBatch = $wb.Worksheets.Item("BatchDecisionOutputValidations")
$wsBatch.Range("F7:F9").ClearContents()

# --- Switch the active tab from BatchDecisionOutputValidations to Cases_RealTimeSpine ---
$wsCases.Activate()
$wsCases.Range("E43").Select()
